# Insert a new data row at row 40 (pushing the existing rows 40-176 down
# to 41-177, matching Excel's native "Insert Row" behaviour) and populate
# the newly created row with the new "Apio" price-record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44707
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = 100112017
$ws.Range("G40").Value = "Apio"
$ws.Range("H40").Value = "Americana (o)"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 700
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 6000
$ws.Range("N40").Value = "`$/docena de matas"
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 1000
$ws.Range("Q40").Value = 6
$ws.Range("R40").Value = "Hortaliza"
